# Auto-generated edit script: applies the scraped diff's cell-value changes
# to the Hades_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 3813
$ws.Range("I94").Value = 1784
$ws.Range("J94").Value = 9900
$ws.Range("K94").Value = 1784
$ws.Range("L94").Value = 9900
$ws.Range("M94").Value = -1333
$ws.Range("N94").Value = -10802
$ws.Range("H107").Value = 996
$ws.Range("I107").Value = 993.6667
$ws.Range("K107").Value = 993.6667
$ws.Range("M107").Value = 926.3333
$ws.Range("H112").Value = 29413698
$ws.Range("I112").Value = 142858100
$ws.Range("J112").Value = 2187.4075
$ws.Range("K112").Value = 428574300
$ws.Range("L112").Value = 6562.2225
$ws.Range("M112").Value = -428573192
$ws.Range("N112").Value = -8778.2225
$ws.Range("H113").Value = 3478.1428
$ws.Range("I113").Value = 3286.3
$ws.Range("J113").Value = 3652.5454
$ws.Range("K113").Value = 3286.3
$ws.Range("L113").Value = 3652.5454
$ws.Range("M113").Value = -32.30000000000018
$ws.Range("N113").Value = -10160.5454
$ws.Range("H137").Value = 1788150
$ws.Range("I137").Value = 2942795.5
$ws.Range("J137").Value = 3697.8635
$ws.Range("K137").Value = 8828386.5
$ws.Range("L137").Value = 11093.5905
$ws.Range("M137").Value = -8825836.5
$ws.Range("N137").Value = -16193.5905
$ws.Range("H138").Value = 3574197.8
$ws.Range("I138").Value = 3032.6667
$ws.Range("J138").Value = 4548152
$ws.Range("K138").Value = 9098.000100000001
$ws.Range("L138").Value = 13644456
$ws.Range("M138").Value = -3958.000100000001
$ws.Range("N138").Value = -13654736
$ws.Range("H140").Value = 74354.25
$ws.Range("J140").Value = 74354.25
$ws.Range("L140").Value = 74354.25
$ws.Range("N140").Value = -84714.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 47716030
$ws.Range("I61").Value = 83418000
$ws.Range("J61").Value = 113410
$ws.Range("K61").Value = 83418000
$ws.Range("L61").Value = 113410
$ws.Range("M61").Value = -83417788
$ws.Range("N61").Value = -113834
$ws.Range("H74").Value = 6845629.5
$ws.Range("I74").Value = 12384163
$ws.Range("J74").Value = 48339.273
$ws.Range("K74").Value = 12384163
$ws.Range("L74").Value = 48339.273
$ws.Range("M74").Value = -12383289
$ws.Range("N74").Value = -50087.273
$ws.Range("H77").Value = 6845629.5
$ws.Range("I77").Value = 12384163
$ws.Range("J77").Value = 48339.273
$ws.Range("K77").Value = 61920815
$ws.Range("L77").Value = 241696.365
$ws.Range("M77").Value = -61916447
$ws.Range("N77").Value = -250432.365
$ws.Range("H122").Value = 13891484
$ws.Range("I122").Value = 1828
$ws.Range("K122").Value = 5484
$ws.Range("M122").Value = -3034
$ws.Range("H132").Value = 34881.812
$ws.Range("I132").Value = 21306.674
$ws.Range("K132").Value = 63920.022
$ws.Range("M132").Value = -61390.022
$ws.Range("H136").Value = 47716030
$ws.Range("I136").Value = 83418000
$ws.Range("J136").Value = 113410
$ws.Range("K136").Value = 250254000
$ws.Range("L136").Value = 340230
$ws.Range("M136").Value = -250251450
$ws.Range("N136").Value = -345330

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 51948.46
$ws.Range("J132").Value = 51948.46
$ws.Range("L132").Value = 51948.46
$ws.Range("N132").Value = -62068.46
$ws.Range("H134").Value = 2052.9048
$ws.Range("I134").Value = 1995.55
$ws.Range("K134").Value = 5986.65
$ws.Range("M134").Value = -3451.65
$ws.Range("H135").Value = 49027.617
$ws.Range("J135").Value = 49027.617
$ws.Range("L135").Value = 49027.617
$ws.Range("N135").Value = -59167.617
$ws.Range("H140").Value = 62937.777
$ws.Range("J140").Value = 62937.777
$ws.Range("L140").Value = 62937.777
$ws.Range("N140").Value = -73297.777
$ws.Range("H141").Value = 66960
$ws.Range("J141").Value = 66960
$ws.Range("L141").Value = 66960
$ws.Range("N141").Value = -77320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""
$ws.Range("H16").Value = 1321.6666
$ws.Range("I16").Value = 1167.7778
$ws.Range("K16").Value = 1167.7778
$ws.Range("M16").Value = -880.7778000000001
$ws.Range("H31").Value = 446243.2
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 446243.2
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 446243.2
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -446833.2
$ws.Range("H34").Value = 446243.2
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 446243.2
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 446243.2
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -446647.2
$ws.Range("H52").Value = 48850
$ws.Range("J52").Value = 48850
$ws.Range("L52").Value = 48850
$ws.Range("N52").Value = -49438
$ws.Range("H86").Value = 3213.5
$ws.Range("I86").Value = 2900
$ws.Range("K86").Value = 2900
$ws.Range("M86").Value = -1777
$ws.Range("H89").Value = 3213.5
$ws.Range("I89").Value = 2900
$ws.Range("K89").Value = 14500
$ws.Range("M89").Value = -8884
$ws.Range("H113").Value = 1321.6666
$ws.Range("I113").Value = 1167.7778
$ws.Range("K113").Value = 1167.7778
$ws.Range("M113").Value = 1002.2222
$ws.Range("H122").Value = 2655.25
$ws.Range("I122").Value = 2625.5715
$ws.Range("J122").Value = 2678.3333
$ws.Range("K122").Value = 7876.7145
$ws.Range("L122").Value = 8034.999899999999
$ws.Range("M122").Value = -5426.7145
$ws.Range("N122").Value = -12934.9999
$ws.Range("H135").Value = 48260
$ws.Range("J135").Value = 48260
$ws.Range("L135").Value = 48260
$ws.Range("N135").Value = -58400
$ws.Range("H139").Value = 52791.668
$ws.Range("J139").Value = 52791.668
$ws.Range("L139").Value = 52791.668
$ws.Range("N139").Value = -63071.668
$ws.Range("H141").Value = 53905.5
$ws.Range("J141").Value = 53905.5
$ws.Range("L141").Value = 53905.5
$ws.Range("N141").Value = -64265.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 607.5238000000001
$ws.Range("I113").Value = 565.8333
$ws.Range("J113").Value = 663.1111
$ws.Range("K113").Value = 1697.4999
$ws.Range("L113").Value = 1989.3333
$ws.Range("M113").Value = 472.5001
$ws.Range("N113").Value = -6329.3333
$ws.Range("H131").Value = 776.1515000000001
$ws.Range("J131").Value = 1059.1428
$ws.Range("L131").Value = 3177.4284
$ws.Range("N131").Value = -13257.4284
$ws.Range("H133").Value = 8700
$ws.Range("I133").Value = 15000
$ws.Range("J133").Value = 8000
$ws.Range("K133").Value = 45000
$ws.Range("L133").Value = 24000
$ws.Range("M133").Value = -39940
$ws.Range("N133").Value = -34120
$ws.Range("H141").Value = 9087.105
$ws.Range("I141").Value = 3237.9167
$ws.Range("J141").Value = 19114.285
$ws.Range("K141").Value = 9713.750100000001
$ws.Range("L141").Value = 57342.855
$ws.Range("M141").Value = -4533.750100000001
$ws.Range("N141").Value = -67702.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1982.8572
$ws.Range("I102").Value = 1966.8235
$ws.Range("J102").Value = 2051
$ws.Range("K102").Value = 1966.8235
$ws.Range("L102").Value = 2051
$ws.Range("M102").Value = -344.8235
$ws.Range("N102").Value = -5295
$ws.Range("H113").Value = 2208.0908
$ws.Range("I113").Value = 1759.8
$ws.Range("J113").Value = 2581.6667
$ws.Range("K113").Value = 1759.8
$ws.Range("L113").Value = 2581.6667
$ws.Range("M113").Value = 410.2
$ws.Range("N113").Value = -6921.6667
$ws.Range("H122").Value = 2619.5454
$ws.Range("I122").Value = 2115.2856
$ws.Range("J122").Value = 3502
$ws.Range("K122").Value = 6345.8568
$ws.Range("L122").Value = 10506
$ws.Range("M122").Value = -3895.8568
$ws.Range("N122").Value = -15406

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 3000
$ws.Range("J11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("N11").Value = -3280
$ws.Range("H16").Value = 1032.5238
$ws.Range("I16").Value = 832.2778
$ws.Range("J16").Value = 2234
$ws.Range("K16").Value = 832.2778
$ws.Range("L16").Value = 2234
$ws.Range("M16").Value = -662.2778
$ws.Range("N16").Value = -2574
$ws.Range("H40").Value = 2608.25
$ws.Range("I40").Value = 2659.9
$ws.Range("J40").Value = 2350
$ws.Range("K40").Value = 2659.9
$ws.Range("L40").Value = 2350
$ws.Range("M40").Value = -2523.9
$ws.Range("N40").Value = -2622
$ws.Range("H132").Value = 48369.137
$ws.Range("I132").Value = 2792.2646
$ws.Range("J132").Value = 203330.5
$ws.Range("K132").Value = 8376.793799999999
$ws.Range("L132").Value = 609991.5
$ws.Range("M132").Value = -5846.793799999999
$ws.Range("N132").Value = -615051.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2650251
$ws.Range("J5").Value = 2650251
$ws.Range("L5").Value = 2650251
$ws.Range("N5").Value = -2650475
$ws.Range("H122").Value = 1133.5834
$ws.Range("I122").Value = 1011.44446
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3034.33338
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -584.33338
$ws.Range("N122").Value = -9400

